$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 255; existing rows 255-275 shift down to 256-276.
$ws.Rows.Item(255).Insert()

$ws.Range("A255").Value2 = 4
$ws.Range("B255").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C255").Value2 = "Los Lagos"
$ws.Range("D255").Value2 = 44769
$ws.Range("E255").Value2 = 10
$ws.Range("F255").Value2 = "Fruta"
$ws.Range("G255").Value2 = 100108
$ws.Range("H255").Value2 = "Tropicales y subtropicales"
$ws.Range("I255").Value2 = 100108005
$ws.Range("J255").Value2 = "Piña"
$ws.Range("K255").Value2 = "Caramelo"
$ws.Range("L255").Value2 = "Primera"
$ws.Range("M255").Value2 = 15
$ws.Range("N255").Value2 = 23000
$ws.Range("O255").Value2 = 23000
$ws.Range("P255").Value2 = 23000
$ws.Range("Q255").Value2 = "$/caja 12 unidades"
$ws.Range("R255").Value2 = "Ecuador"
$ws.Range("S255").Value2 = 1917
$ws.Range("T255").Value2 = 12
